$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"

# Column D carries the date-formatted style already used by the other rows.
# Use the raw serial date number (44516 = 2021-11-16) to avoid a time-of-day component.
$ws.Cells.Item($row, 4).Value = 44516
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = "Arándano (blue)"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 80
$ws.Cells.Item($row, 14).Value = 3700
$ws.Cells.Item($row, 15).Value = 3800
$ws.Cells.Item($row, 16).Value = 3750
$ws.Cells.Item($row, 17).Value = "$/kilo"
$ws.Cells.Item($row, 18).Value = "Región del Maule"
$ws.Cells.Item($row, 19).Value = 3750
$ws.Cells.Item($row, 20).Value = 1
